# Update "paises.xlsx" country/provincia data per the 08:22 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 08:22"

# --- Country name re-orderings (rows keep their position, but the
#     country label attached to the row changes because the underlying
#     shared-string table was re-sorted) ---

# Kazajistan now ranks above Tailandia
$ws.Cells.Item(61, 1).Value = "Kazajistan"
$ws.Cells.Item(62, 1).Value = "Tailandia"

# Honduras now ranks above Niger / Costa Rica / Kirguistan (cascades down one row)
$ws.Cells.Item(99, 1).Value  = "Honduras"
$ws.Cells.Item(100, 1).Value = "Niger"
$ws.Cells.Item(101, 1).Value = "Costa Rica"
$ws.Cells.Item(102, 1).Value = "Kirguistan"

# Guayana Francesa now ranks above Camboya / Trinidad yTobago (cascades down one row)
$ws.Cells.Item(143, 1).Value = "Guayana Francesa"
$ws.Cells.Item(144, 1).Value = "Camboya"
$ws.Cells.Item(145, 1).Value = "Trinidad yTobago"

# --- Updated case counts (B=Casos totales, C=Nuevos casos, D=Casos activos,
#     E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 34 - Rumania
Set-Row 34 11339 0 3141 7552 227 5 646

# Row 61 - Kazajistan (fresh numbers)
Set-Row 61 2950 115 725 2200 41 0 25

# Row 62 - Tailandia (carries the old row-61 numbers forward)
Set-Row 62 2938 7 2652 232 61 2 54

# Row 64 - Hungria
Set-Row 64 2649 66 516 1842 49 11 291

# Row 99 - Honduras (fresh numbers)
Set-Row 99 702 41 79 559 10 3 64

# Row 100 - Niger (carries the old row-99 numbers forward)
Set-Row 100 701 0 385 287 0 0 29

# Row 101 - Costa Rica (carries the old row-100 numbers forward)
Set-Row 101 697 0 287 404 8 0 6

# Row 102 - Kirguistan (carries the old row-101 numbers forward)
Set-Row 102 695 0 395 292 13 0 8

# Row 108 - Georgia
Set-Row 108 511 14 156 349 6 0 6

# Row 109 - Somalia
Set-Row 109 480 0 14 440 2 0 26

# Row 113 - Taiwan
Set-Row 113 429 0 307 116 0 0 6

# Row 143 - Guayana Francesa (fresh numbers)
Set-Row 143 124 13 91 32 0 0 1

# Row 144 - Camboya (carries the old row-143 numbers forward)
Set-Row 144 122 0 119 3 1 0 0

# Row 145 - Trinidad yTobago (carries the old row-144 numbers forward)
Set-Row 145 116 0 59 49 0 0 8

# Row 162 - Islas Caimanes
Set-Row 162 70 0 10 59 3 0 1
